$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 257, shifting existing rows 257:263 down to 258:264.
$ws.Rows.Item(257).Insert()

# Populate the new row 257 with a new Membrillo price record, matching the
# constant columns of the surrounding records and new values for the rest.
$ws.Cells.Item(257, 1).Value = 10
$ws.Cells.Item(257, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(257, 3).Value = "La Araucanía"
$ws.Cells.Item(257, 4).Value = 45041
$ws.Cells.Item(257, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(257, 5).Value = 9
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100104
$ws.Cells.Item(257, 8).Value = "Frutos de pepita"
$ws.Cells.Item(257, 9).Value = 100104003
$ws.Cells.Item(257, 10).Value = "Membrillo"
$ws.Cells.Item(257, 11).Value = "Champion"
$ws.Cells.Item(257, 12).Value = "Primera"
$ws.Cells.Item(257, 13).Value = 65
$ws.Cells.Item(257, 14).Value = 13000
$ws.Cells.Item(257, 15).Value = 13000
$ws.Cells.Item(257, 16).Value = 13000
$ws.Cells.Item(257, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(257, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(257, 19).Value = 722
$ws.Cells.Item(257, 20).Value = 18
